# Update crypto price (D) and 1h volume/change (E) columns per the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.873.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.03%  "

$ws.Range("D3").Value = "'3.466.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.00%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'604.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.16%  "

$ws.Range("D6").Value = "'148.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.73%  "

$ws.Range("D7").Value = "'3.462.67"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.09%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -2.20%  "

$ws.Range("E10").Value = "  -4.71%  "

$ws.Range("D11").Value = "'7.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("D12").Value = "'0.425"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.94%  "

$ws.Range("E13").Value = "  -4.99%  "

$ws.Range("D14").Value = "'31.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.72%  "

$ws.Range("D15").Value = "'4.047.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.16%  "

$ws.Range("D16").Value = "'3.459.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.26%  "

$ws.Range("D17").Value = "'66.878.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.82%  "

$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").Value = "'6.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.26%  "

$ws.Range("D20").Value = "'15.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.66%  "

$ws.Range("E21").Value = "  -2.35%  "

$ws.Range("D22").Value = "'440.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.97%  "

$ws.Range("D23").Value = "'0.611"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.66%  "

$ws.Range("D24").Value = "'78.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.34%  "

$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").Value = "'3.599.23"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.20%  "

$ws.Range("E27").Value = "  -10.28%  "

$ws.Range("D28").Value = "'9.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.41%  "

$ws.Range("D29").Value = "'8.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -10.07%  "

$ws.Range("E30").Value = "  -6.47%  "

$ws.Range("E31").Value = "  -7.12%  "

$ws.Range("E32").Value = "  -3.10%  "

$ws.Range("E33").Value = "  +0.04%  "

$ws.Range("D34").Value = "'25.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.29%  "

$ws.Range("E35").Value = "  -7.57%  "

$ws.Range("D36").Value = "'3.454.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.28%  "

$ws.Range("E37").Value = "  -7.83%  "

$ws.Range("E38").Value = "  -6.59%  "

$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "'173.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.06%  "

$ws.Range("D42").Value = "'0.0893"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.79%  "

$ws.Range("E43").Value = "  -11.36%  "

$ws.Range("E44").Value = "  -5.06%  "

$ws.Range("D45").Value = "'0.885"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.23%  "

$ws.Range("D46").Value = "'29.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.92%  "

$ws.Range("E47").Value = "  -0.23%  "

$ws.Range("E48").Value = "  -10.47%  "

$ws.Range("E49").Value = "  -10.89%  "

$ws.Range("D50").Value = "'7.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.41%  "

$ws.Range("D51").Value = "'0.987"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.60%  "
